$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at D:E, shifting existing D:M (old D:K) to F:M
$ws.Range("D:E").Insert()

# Copy number formats from column F (the old column D, now shifted right by 2)
# onto the two new columns D:E, row-range by row-range, so that we only touch
# rows that actually have formatted cells (skip section-header rows that have
# no D..M content at all).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New values for the two inserted columns (D, E) per row.
$newValues = @{
    7 = @(43465, 43373)
    8 = @(924200, 1016200)
    9 = @(386700, 422500)
    10 = @(537500, 593700)
    11 = @($null, $null)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(0, 14100)
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(617700, 676600)
    18 = @(306500, 339600)
    19 = @($null, $null)
    20 = @(4400, 2900)
    21 = @(325400, 356800)
    22 = @(0, 0)
    23 = @(310900, 342600)
    24 = @(70000, 74800)
    25 = @(0, 0)
    26 = @(240900, 267700)
    27 = @(240900, 267700)
    28 = @(0, 0)
    29 = @(-1800, "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-4400, -2900)
    33 = @(239100, 267700)
    34 = @(0, 0)
    35 = @(239100, 267700)
    38 = @(43465, 43373)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(637500, 713700)
    42 = @(320700, 457900)
    43 = @(484600, 620200)
    44 = @(277700, 262100)
    45 = @(83700, 98800)
    46 = @(1804200, 2152700)
    47 = @(0, 1600)
    48 = @(243100, 242900)
    49 = @(2377500, 2373900)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(102100, 100300)
    53 = @(0, 0)
    54 = @(4526900, 4871400)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(248800, 278900)
    58 = @(0, 0)
    59 = @(352400, 377400)
    60 = @(601100, 656400)
    61 = @(0, 0)
    62 = @(314800, 321700)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(916000, 978100)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(3914600, 3675500)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(3610900, 3893300)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(239100, 267700)
    82 = @($null, $null)
    83 = @(14500, 14300)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(340500, 320100)
    90 = @($null, $null)
    91 = @(-12100, -15200)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(118800, -272000)
    95 = @($null, $null)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-534000, 10400)
    101 = @(-1600, -4400)
    102 = @(-76200, 54000)
}

foreach ($r in $newValues.Keys) {
    $vals = $newValues[$r]
    if ($vals[0] -ne $null) { $ws.Cells.Item($r, 4).Value = $vals[0] }
    if ($vals[1] -ne $null) { $ws.Cells.Item($r, 5).Value = $vals[1] }
}
